$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 79: 2024-01-20 13:02:43 resale snapshot
$ws.Range("A79").NumberFormat = "@"
$ws.Range("A79").Value = "2024-01-20"
$ws.Range("A79").ClearFormats()

$ws.Range("B79").NumberFormat = "@"
$ws.Range("B79").Value = "13:02:43"
$ws.Range("B79").ClearFormats()

$ws.Range("C79").Value = "Saturday"

$ws.Range("D79").NumberFormat = "@"
$ws.Range("D79").Value = "02"
$ws.Range("D79").ClearFormats()

$ws.Range("E79").Value = 138585
$ws.Range("F79").Value = 140788
$ws.Range("G79").Value = 171481
$ws.Range("H79").Value = 148906
$ws.Range("I79").Value = -1
$ws.Range("J79").Value = 122576
$ws.Range("K79").Value = 223642
$ws.Range("L79").Value = 255247
$ws.Range("M79").Value = 185343
$ws.Range("N79").Value = 110385
$ws.Range("O79").Value = 41158
$ws.Range("P79").Value = 30919
$ws.Range("Q79").Value = 73602
$ws.Range("R79").Value = -1
$ws.Range("S79").Value = 42626
$ws.Range("T79").Value = -1
